$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''67.025.75'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.71%  '

# Row 3
$ws.Range('D3').Value = '''3.254.53'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.68%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').Value = '''578.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.72%  '

# Row 6
$ws.Range('D6').Value = '''177.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.42%  '

# Row 7
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$ws.Range('D8').Value = '''0.600'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.48%  '

# Row 9
$ws.Range('D9').Value = '''3.252.83'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.85%  '

# Row 10
$ws.Range('E10').Value = '  +4.32%  '

# Row 11
$ws.Range('D11').Value = '''6.74'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.81%  '

# Row 12
$ws.Range('E12').Value = '  +4.10%  '

# Row 13
$ws.Range('D13').Value = '''3.822.50'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.92%  '

# Row 14
$ws.Range('E14').Value = '  +0.63%  '

# Row 15
$ws.Range('D15').Value = '''28.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.79%  '

# Row 16
$ws.Range('D16').Value = '''66.997.08'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.65%  '

# Row 17
$ws.Range('E17').Value = '  +3.20%  '

# Row 18
$ws.Range('D18').Value = '''3.252.94'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.69%  '

# Row 19
$ws.Range('D19').Value = '''5.85'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.51%  '

# Row 20
$ws.Range('E20').Value = '  +2.72%  '

# Row 21
$ws.Range('D21').Value = '''372.93'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.58%  '

# Row 22
$ws.Range('E22').Value = '  +5.78%  '

# Row 23
$ws.Range('E23').Value = '  +0.27%  '

# Row 24
$ws.Range('D24').Value = '''71.17'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.12%  '

# Row 25
$ws.Range('D25').Value = '''0.512'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.03%  '

# Row 26
$ws.Range('D26').Value = '''3.396.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.66%  '

# Row 27
$ws.Range('D27').Value = '''0.0000119'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.98%  '

# Row 28
$ws.Range('D28').Value = '''9.83'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.57%  '

# Row 29
$ws.Range('D29').Value = '''0.179'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.29%  '

# Row 30
$ws.Range('D30').Value = '''0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.20%  '

# Row 31
$ws.Range('E31').Value = '  +4.15%  '

# Row 32
$ws.Range('D32').Value = '''5.60'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.60%  '

# Row 33
$ws.Range('D33').Value = '''22.60'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.35%  '

# Row 34
$ws.Range('E34').Value = '  +0.05%  '

# Row 35
$ws.Range('E35').Value = '  +5.07%  '

# Row 36
$ws.Range('D36').Value = '''6.82'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.05%  '

# Row 37
$ws.Range('D37').Value = '''167.25'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.80%  '

# Row 38
$ws.Range('E38').Value = '  +4.51%  '

# Row 39
$ws.Range('D39').Value = '''0.853'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.56%  '

# Row 40
$ws.Range('D40').Value = '''1.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.91%  '

# Row 41
$ws.Range('D41').Value = '''27.20'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.36%  '

# Row 42
$ws.Range('E42').Value = '  +1.52%  '

# Row 43
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '''2.760.13'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.05%  '

# Row 44
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''6.45'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +7.14%  '

# Row 45
$ws.Range('E45').Value = '  +4.53%  '

# Row 46
$ws.Range('D46').Value = '''349.63'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.22%  '

# Row 47
$ws.Range('D47').Value = '''25.27'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.69%  '

# Row 48
$ws.Range('D48').Value = '''40.50'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.66%  '

# Row 49
$ws.Range('E49').Value = '  +2.76%  '

# Row 50
$ws.Range('E50').Value = '  +3.15%  '

# Row 51
$ws.Range('E51').Value = '  +1.04%  '
